$wb = $excel.ActiveWorkbook

$newGuid = "dea3305b-635d-45fa-af3a-19e14bc3d44d"
$newHoHash = "fdb8bbdb6ef1fc41ac4ed20b851a005a639d2c69"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-24 22:58:07"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.$newHoHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-24 22:57:57"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.$newHoHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-24 22:58:07"
